# chore: update Sheets via scheduled runner
# Refresh cached marketboard price/profit figures (columns H-N) for the
# affected leve rows across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 295.81818
$ws.Range("I53").Value = 89.375
$ws.Range("J53").Value = 413.7857
$ws.Range("K53").Value = 89.375
$ws.Range("L53").Value = 413.7857
$ws.Range("M53").Value = 547.625
$ws.Range("N53").Value = -1687.7857

$ws.Range("H99").Value = 8793.857
$ws.Range("I99").Value = 259.5
$ws.Range("J99").Value = 60000
$ws.Range("K99").Value = 778.5
$ws.Range("L99").Value = 180000
$ws.Range("M99").Value = 719.5
$ws.Range("N99").Value = -182996

$ws.Range("H138").Value = 3766.1755
$ws.Range("J138").Value = 4558.7046
$ws.Range("L138").Value = 13676.1138
$ws.Range("N138").Value = -23956.1138

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 2001902
$ws.Range("I11").Value = 2501376.5
$ws.Range("J11").Value = 4004
$ws.Range("K11").Value = 2501376.5
$ws.Range("L11").Value = 4004
$ws.Range("M11").Value = -2501232.5
$ws.Range("N11").Value = -4292

$ws.Range("H24").Value = 36416.668
$ws.Range("J24").Value = 36416.668
$ws.Range("L24").Value = 36416.668
$ws.Range("N24").Value = -37164.668

$ws.Range("H61").Value = 2837.818
$ws.Range("I61").Value = 2230.8572
$ws.Range("J61").Value = 3900
$ws.Range("K61").Value = 2230.8572
$ws.Range("L61").Value = 3900
$ws.Range("M61").Value = -2018.8572
$ws.Range("N61").Value = -4324

$ws.Range("H74").Value = 2739.0698
$ws.Range("I74").Value = 2353.4614
$ws.Range("K74").Value = 2353.4614
$ws.Range("M74").Value = -1479.4614

$ws.Range("H77").Value = 2739.0698
$ws.Range("I77").Value = 2353.4614
$ws.Range("K77").Value = 11767.307
$ws.Range("M77").Value = -7399.307000000001

$ws.Range("H100").Value = 36416.668
$ws.Range("J100").Value = 36416.668
$ws.Range("L100").Value = 36416.668
$ws.Range("N100").Value = -38580.668

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H122").Value = 1825.6
$ws.Range("I122").Value = 1529.8948
$ws.Range("K122").Value = 4589.6844
$ws.Range("M122").Value = -2139.6844

$ws.Range("H132").Value = 1856.8889
$ws.Range("I132").Value = 1788.8235
$ws.Range("K132").Value = 5366.470499999999
$ws.Range("M132").Value = -2836.470499999999

$ws.Range("H136").Value = 2837.818
$ws.Range("I136").Value = 2230.8572
$ws.Range("J136").Value = 3900
$ws.Range("K136").Value = 6692.571599999999
$ws.Range("L136").Value = 11700
$ws.Range("M136").Value = -4142.571599999999
$ws.Range("N136").Value = -16800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2376.4285
$ws.Range("I134").Value = 2021.7307
$ws.Range("K134").Value = 6065.1921
$ws.Range("M134").Value = -3530.1921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 267.16666
$ws.Range("I7").Value = 85.5
$ws.Range("K7").Value = 85.5
$ws.Range("M7").Value = 27.5

$ws.Range("H58").Value = 2471.3157
$ws.Range("I58").Value = 1531.3334
$ws.Range("J58").Value = 5996.25
$ws.Range("K58").Value = 1531.3334
$ws.Range("L58").Value = 5996.25
$ws.Range("M58").Value = -1328.3334
$ws.Range("N58").Value = -6402.25

$ws.Range("H69").Value = 19087.857
$ws.Range("I69").Value = 10573
$ws.Range("K69").Value = 10573
$ws.Range("M69").Value = -9824

$ws.Range("H72").Value = 19087.857
$ws.Range("I72").Value = 10573
$ws.Range("K72").Value = 31719
$ws.Range("M72").Value = -27975

$ws.Range("H107").Value = 10241.619
$ws.Range("I107").Value = 569.0769
$ws.Range("K107").Value = 569.0769
$ws.Range("M107").Value = 1350.9231

$ws.Range("H132").Value = 1920.8462
$ws.Range("I132").Value = 1890.9565
$ws.Range("K132").Value = 5672.8695
$ws.Range("M132").Value = -3142.8695

$ws.Range("H134").Value = 2045.8636
$ws.Range("I134").Value = 2022.3903
$ws.Range("K134").Value = 6067.1709
$ws.Range("M134").Value = -3532.1709

$ws.Range("H136").Value = 2471.3157
$ws.Range("I136").Value = 1531.3334
$ws.Range("J136").Value = 5996.25
$ws.Range("K136").Value = 4594.0002
$ws.Range("L136").Value = 17988.75
$ws.Range("M136").Value = -2044.0002
$ws.Range("N136").Value = -23088.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 89973
$ws.Range("I80").Value = 164667.28
$ws.Range("J80").Value = 2829.6667
$ws.Range("K80").Value = 164667.28
$ws.Range("L80").Value = 2829.6667
$ws.Range("M80").Value = -163669.28
$ws.Range("N80").Value = -4825.6667

$ws.Range("H83").Value = 89973
$ws.Range("I83").Value = 164667.28
$ws.Range("J83").Value = 2829.6667
$ws.Range("K83").Value = 823336.4
$ws.Range("L83").Value = 14148.3335
$ws.Range("M83").Value = -818344.4
$ws.Range("N83").Value = -24132.3335

$ws.Range("H132").Value = 2784.3438
$ws.Range("I132").Value = 2362.5217
$ws.Range("J132").Value = 3862.3333
$ws.Range("K132").Value = 7087.5651
$ws.Range("L132").Value = 11586.9999
$ws.Range("M132").Value = -4557.5651
$ws.Range("N132").Value = -16646.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2032.4445
$ws.Range("I16").Value = 2484.7144
$ws.Range("J16").Value = 449.5
$ws.Range("K16").Value = 2484.7144
$ws.Range("L16").Value = 449.5
$ws.Range("M16").Value = -2314.7144
$ws.Range("N16").Value = -789.5

$ws.Range("H46").Value = 2657.4546
$ws.Range("J46").Value = 2777.0527
$ws.Range("L46").Value = 2777.0527
$ws.Range("N46").Value = -3153.0527

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H132").Value = 3818.1777
$ws.Range("I132").Value = 2456.6155
$ws.Range("J132").Value = 12668.333
$ws.Range("K132").Value = 7369.8465
$ws.Range("L132").Value = 38004.999
$ws.Range("M132").Value = -4839.8465
$ws.Range("N132").Value = -43064.999

$ws.Range("H136").Value = 2834.0613
$ws.Range("J136").Value = 2554.0908
$ws.Range("L136").Value = 7662.2724
$ws.Range("N136").Value = -12762.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1000000
$ws.Range("I8").Value = 1000000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1000000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -999860
$ws.Range("N8").ClearContents()

$ws.Range("H81").Value = 3814.7827
$ws.Range("I81").Value = 3235.6667
$ws.Range("K81").Value = 6471.3334
$ws.Range("M81").Value = -5410.3334

$ws.Range("H84").Value = 3814.7827
$ws.Range("I84").Value = 3235.6667
$ws.Range("K84").Value = 32356.667
$ws.Range("M84").Value = -27052.667

$ws.Range("H107").Value = 593.1177
$ws.Range("I107").Value = 659.8570999999999
$ws.Range("J107").Value = 281.66666
$ws.Range("K107").Value = 1979.5713
$ws.Range("L107").Value = 844.9999799999999
$ws.Range("M107").Value = -59.57129999999984
$ws.Range("N107").Value = -4684.99998

$ws.Range("H122").Value = 3867.9285
$ws.Range("I122").Value = 3753.7273
$ws.Range("K122").Value = 11261.1819
$ws.Range("M122").Value = -8811.1819
